$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit "Adding the changes we made on may 9th": the accelerometer/gyro readings (ax,ay,az,gx,gy,gz)
# for the existing "falling" sample window were re-sliced 8 rows later in the source capture,
# and the window was extended with 8 new leading samples and 2 new trailing samples.
# The timestamp (A) and label (B) columns stay a plain sequential 0,100,200,... "falling" run
# across the whole (now longer) range; only the sensor columns C:H were overwritten/extended.

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -1.275631546974182
$ws.Cells.Item(2, 4).Value = 1.748281717300415
$ws.Cells.Item(2, 5).Value = 0.7527783811092381
$ws.Cells.Item(2, 6).Value = -0.0042760567739605
$ws.Cells.Item(2, 7).Value = -0.1111774742603302
$ws.Cells.Item(2, 8).Value = -0.0980438739061355

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -1.113769233226775
$ws.Cells.Item(3, 4).Value = 1.769958406686783
$ws.Cells.Item(3, 5).Value = 1.082688376307487
$ws.Cells.Item(3, 6).Value = 0.0734565481543541
$ws.Cells.Item(3, 7).Value = 0.1905899494886398
$ws.Cells.Item(3, 8).Value = -0.1611157059669494

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = 0.6056947708129873
$ws.Cells.Item(4, 4).Value = 1.413846492767334
$ws.Cells.Item(4, 5).Value = 1.041245818138122
$ws.Cells.Item(4, 6).Value = -0.101709060370922
$ws.Cells.Item(4, 7).Value = -0.09025534242391579
$ws.Cells.Item(4, 8).Value = 0.1893682330846786

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -0.2511940002441411
$ws.Cells.Item(5, 4).Value = 1.83873063325882
$ws.Cells.Item(5, 5).Value = 0.5010688602924345
$ws.Cells.Item(5, 6).Value = -0.0087048299610614
$ws.Cells.Item(5, 7).Value = 0.024892758578062
$ws.Cells.Item(5, 8).Value = 0.00167987938039

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -0.7442607879638676
$ws.Cells.Item(6, 4).Value = 1.761505782604217
$ws.Cells.Item(6, 5).Value = 0.9367214739322666
$ws.Cells.Item(6, 6).Value = 0.0125227374956011
$ws.Cells.Item(6, 7).Value = -0.0652098655700683
$ws.Cells.Item(6, 8).Value = 0.066737025976181

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = -0.8476336002349854
$ws.Cells.Item(7, 4).Value = 1.69824892282486
$ws.Cells.Item(7, 5).Value = 0.9451412782073016
$ws.Cells.Item(7, 6).Value = -0.0397062413394451
$ws.Cells.Item(7, 7).Value = 0.0247400421649217
$ws.Cells.Item(7, 8).Value = 0.030695978552103

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = -0.6413483619689934
$ws.Cells.Item(8, 4).Value = 1.650843858718872
$ws.Cells.Item(8, 5).Value = 0.9322790801525122
$ws.Cells.Item(8, 6).Value = 0.012980886735022
$ws.Cells.Item(8, 7).Value = 0.0633772686123848
$ws.Cells.Item(8, 8).Value = -0.0366519130766391

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "falling"
$ws.Cells.Item(9, 3).Value = -0.5721501111984255
$ws.Cells.Item(9, 4).Value = 1.609763711690903
$ws.Cells.Item(9, 5).Value = 1.015784159302711
$ws.Cells.Item(9, 6).Value = -0.0360410511493682
$ws.Cells.Item(9, 7).Value = 0.0274889357388019
$ws.Cells.Item(9, 8).Value = -0.0474947728216648

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "falling"
$ws.Cells.Item(10, 3).Value = -0.6957695484161374
$ws.Cells.Item(10, 4).Value = 1.588029444217682
$ws.Cells.Item(10, 5).Value = 1.020436197519302
$ws.Cells.Item(10, 6).Value = -0.0114537235349416
$ws.Cells.Item(10, 7).Value = -0.0215329993516206
$ws.Cells.Item(10, 8).Value = 0.0035124751739203

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "falling"
$ws.Cells.Item(11, 3).Value = -0.4331760406494141
$ws.Cells.Item(11, 4).Value = 1.524660766124726
$ws.Cells.Item(11, 5).Value = 1.148605212569236
$ws.Cells.Item(11, 6).Value = -0.0171042270958423
$ws.Cells.Item(11, 7).Value = -0.0526871271431446
$ws.Cells.Item(11, 8).Value = 0.0519235469400882

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "falling"
$ws.Cells.Item(12, 3).Value = -0.3733081817626966
$ws.Cells.Item(12, 4).Value = 1.395935773849488
$ws.Cells.Item(12, 5).Value = 1.908020853996279
$ws.Cells.Item(12, 6).Value = 0.0067195175215601
$ws.Cells.Item(12, 7).Value = 0.1915062516927719
$ws.Cells.Item(12, 8).Value = -0.0238237436860799

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "falling"
$ws.Cells.Item(13, 3).Value = -0.9485191106796277
$ws.Cells.Item(13, 4).Value = 1.76411008834839
$ws.Cells.Item(13, 5).Value = 2.944827482104304
$ws.Cells.Item(13, 6).Value = -0.001527163083665
$ws.Cells.Item(13, 7).Value = 0.3747658133506775
$ws.Cells.Item(13, 8).Value = 0.0192422550171613

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "falling"
$ws.Cells.Item(14, 3).Value = -1.393463373184205
$ws.Cells.Item(14, 4).Value = 1.543532192707062
$ws.Cells.Item(14, 5).Value = 3.777262568473817
$ws.Cells.Item(14, 6).Value = -0.3475823104381561
$ws.Cells.Item(14, 7).Value = 0.5015203952789307
$ws.Cells.Item(14, 8).Value = -0.2428189367055893

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "falling"
$ws.Cells.Item(15, 3).Value = -1.7527374625206
$ws.Cells.Item(15, 4).Value = 0.9714637398719776
$ws.Cells.Item(15, 5).Value = 4.131629109382629
$ws.Cells.Item(15, 6).Value = -0.2506074607372284
$ws.Cells.Item(15, 7).Value = 0.8736900091171265
$ws.Cells.Item(15, 8).Value = -0.1963931769132614

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "falling"
$ws.Cells.Item(16, 3).Value = -2.943279504776003
$ws.Cells.Item(16, 4).Value = 0.2452936768531791
$ws.Cells.Item(16, 5).Value = 4.232949018478394
$ws.Cells.Item(16, 6).Value = -0.168751522898674
$ws.Cells.Item(16, 7).Value = 0.936914563179016
$ws.Cells.Item(16, 8).Value = -0.131183311343193

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "falling"
$ws.Cells.Item(17, 3).Value = -3.76435226202011
$ws.Cells.Item(17, 4).Value = 0.9894824773073276
$ws.Cells.Item(17, 5).Value = 3.341848820447916
$ws.Cells.Item(17, 6).Value = -0.2408336251974105
$ws.Cells.Item(17, 7).Value = -0.064446285367012
$ws.Cells.Item(17, 8).Value = -0.3248275816440582

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "falling"
$ws.Cells.Item(18, 3).Value = -0.01425111293788817
$ws.Cells.Item(18, 4).Value = 1.767028868198381
$ws.Cells.Item(18, 5).Value = -1.569932878017464
$ws.Cells.Item(18, 6).Value = -0.2364048510789871
$ws.Cells.Item(18, 7).Value = 0.18539759516716
$ws.Cells.Item(18, 8).Value = -1.519680023193359

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "falling"
$ws.Cells.Item(19, 3).Value = 7.34391605854027
$ws.Cells.Item(19, 4).Value = -0.5457401573657639
$ws.Cells.Item(19, 5).Value = -8.302204966545009
$ws.Cells.Item(19, 6).Value = -1.216843605041504
$ws.Cells.Item(19, 7).Value = -1.400255799293518
$ws.Cells.Item(19, 8).Value = -1.343445420265198

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "falling"
$ws.Cells.Item(20, 3).Value = -11.97496986389153
$ws.Cells.Item(20, 4).Value = 8.297744750976534
$ws.Cells.Item(20, 5).Value = 15.6360607147216
$ws.Cells.Item(20, 6).Value = 1.094823241233826
$ws.Cells.Item(20, 7).Value = -2.151772737503052
$ws.Cells.Item(20, 8).Value = 0.7765624523162842

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "falling"
$ws.Cells.Item(21, 3).Value = 3.999401926994334
$ws.Cells.Item(21, 4).Value = 1.987441882491108
$ws.Cells.Item(21, 5).Value = -1.899090290069591
$ws.Cells.Item(21, 6).Value = -0.2055561542510986
$ws.Cells.Item(21, 7).Value = -1.607033729553223
$ws.Cells.Item(21, 8).Value = -0.9816604256629944

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "falling"
$ws.Cells.Item(22, 3).Value = 2.015434503555293
$ws.Cells.Item(22, 4).Value = 2.558928638696672
$ws.Cells.Item(22, 5).Value = -0.2935100495815224
$ws.Cells.Item(22, 6).Value = -0.5111414790153503
$ws.Cells.Item(22, 7).Value = -1.858557462692261
$ws.Cells.Item(22, 8).Value = 0.0117591563612222

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "falling"
$ws.Cells.Item(23, 3).Value = 0.5337359905242919
$ws.Cells.Item(23, 4).Value = 2.63483315706253
$ws.Cells.Item(23, 5).Value = -0.1261084899306326
$ws.Cells.Item(23, 6).Value = -0.3014619946479797
$ws.Cells.Item(23, 7).Value = -0.9010262489318848
$ws.Cells.Item(23, 8).Value = 0.4983133375644684

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "falling"
$ws.Cells.Item(24, 3).Value = 1.104637145996096
$ws.Cells.Item(24, 4).Value = 1.630028188228604
$ws.Cells.Item(24, 5).Value = -0.8494508564472198
$ws.Cells.Item(24, 6).Value = -0.2797762751579284
$ws.Cells.Item(24, 7).Value = -1.265254616737366
$ws.Cells.Item(24, 8).Value = 0.2492330223321914

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "falling"
$ws.Cells.Item(25, 3).Value = -1.313743948936482
$ws.Cells.Item(25, 4).Value = 1.221044480800631
$ws.Cells.Item(25, 5).Value = -1.651125282049184
$ws.Cells.Item(25, 6).Value = 0.2038762718439102
$ws.Cells.Item(25, 7).Value = -1.481042742729187
$ws.Cells.Item(25, 8).Value = -0.2376265823841095

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "falling"
$ws.Cells.Item(26, 3).Value = -5.144322872161865
$ws.Cells.Item(26, 4).Value = 1.553820222616196
$ws.Cells.Item(26, 5).Value = -2.058750659227372
$ws.Cells.Item(26, 6).Value = -0.328340083360672
$ws.Cells.Item(26, 7).Value = 0.6785185933113098
$ws.Cells.Item(26, 8).Value = 0.7831292152404785

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "falling"
$ws.Cells.Item(27, 3).Value = -1.258464217185967
$ws.Cells.Item(27, 4).Value = 1.119151741266259
$ws.Cells.Item(27, 5).Value = 0.6321565061807609
$ws.Cells.Item(27, 6).Value = 0.0797179117798805
$ws.Cells.Item(27, 7).Value = 0.1185078546404838
$ws.Cells.Item(27, 8).Value = 0.2952006161212921

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "falling"
$ws.Cells.Item(28, 3).Value = 0.5062999725341749
$ws.Cells.Item(28, 4).Value = 3.225210666656487
$ws.Cells.Item(28, 5).Value = 0.02823758125305291
$ws.Cells.Item(28, 6).Value = 0.0479529201984405
$ws.Cells.Item(28, 7).Value = 0.3949243724346161
$ws.Cells.Item(28, 8).Value = -0.1617265790700912

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "falling"
$ws.Cells.Item(29, 3).Value = -0.5733692646026665
$ws.Cells.Item(29, 4).Value = 1.683094680309288
$ws.Cells.Item(29, 5).Value = 0.3010409921407713
$ws.Cells.Item(29, 6).Value = -0.08704829961061469
$ws.Cells.Item(29, 7).Value = -0.058643065392971
$ws.Cells.Item(29, 8).Value = -0.2785545587539673

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "falling"
$ws.Cells.Item(30, 3).Value = -0.4257340431213359
$ws.Cells.Item(30, 4).Value = 1.845006287097933
$ws.Cells.Item(30, 5).Value = 0.8945446908474008
$ws.Cells.Item(30, 6).Value = -0.0519235469400882
$ws.Cells.Item(30, 7).Value = -0.107512280344963
$ws.Cells.Item(30, 8).Value = 0.1065959855914115

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "falling"
$ws.Cells.Item(31, 3).Value = 0.06665813922882169
$ws.Cells.Item(31, 4).Value = 1.836877554655074
$ws.Cells.Item(31, 5).Value = 0.7217497229576104
$ws.Cells.Item(31, 6).Value = 0.0108428578823804
$ws.Cells.Item(31, 7).Value = -0.0210748501121997
$ws.Cells.Item(31, 8).Value = -0.0189368221908807
